# "User Administrator Account added and working"
#
# The document currently starts with:
#   1) "User Administrator Account" paragraph (carries the hidden
#      "_GoBack" bookmark)
#   2) two empty placeholder paragraphs
#   3) "Search Function (All Users)" paragraph (the heading that should
#      own the "_GoBack" bookmark going forward)
#
# The edit removes paragraphs 1-3 and relocates the "_GoBack" bookmark to
# the very start of the (now first) "Search Function (All Users)"
# paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that will become the new first paragraph
# (paragraph 4: "Search Function (All Users)").
$target = $d.Paragraphs.Item(4)

# Re-create the "_GoBack" bookmark at the start of that paragraph while it
# still has a paragraph before it (placing a zero-length bookmark at the
# very first character position of the document snaps its end to the
# following paragraph, so do this before the preceding paragraphs are
# removed).
$startPos = $target.Range.Start
$bookmarkRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the leading "User Administrator Account" paragraph plus the two
# empty paragraphs that followed it.
$removeRange = $d.Range($d.Paragraphs.Item(1).Range.Start, $d.Paragraphs.Item(3).Range.End)
$removeRange.Delete()
